# Restore cell C10 on the "Rules" sheet back to its original value of 1
# (it had been changed to 18).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
